# Update column F ("想去人数") values on several rows across the four
# worksheets of the workbook, incrementing the "want to go" counters to
# match the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> list of (row, newValue) pairs for column F.
$changes = @{
    "展览"     = @{ 14 = 631; 16 = 1360; 17 = 19; 19 = 3809; 22 = 757; 23 = 22; 32 = 1039; 33 = 1033 }
    "演出"     = @{ 20 = 477 }
    "本地生活" = @{ 4 = 530 }
    "全部类型" = @{ 9 = 530; 27 = 1360; 28 = 19; 31 = 3809; 34 = 757; 45 = 477; 48 = 1039; 49 = 1033 }
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowsForSheet = $changes[$sheetName]
    foreach ($row in $rowsForSheet.Keys) {
        $newValue = $rowsForSheet[$row]
        $ws.Range("F$row").Value = $newValue
    }
}
